$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(15, 8).Value = 1794549.9  # H15: 1823965.9 -> 1794549.9
$ws.Cells.Item(15, 9).Value = 1794549.9  # I15: 1823965.9 -> 1794549.9
$ws.Cells.Item(15, 11).Value = 5383649.699999999  # K15: 5471897.699999999 -> 5383649.699999999
$ws.Cells.Item(15, 13).Value = -5383480.699999999  # M15: -5471728.699999999 -> -5383480.699999999

$ws.Cells.Item(33, 8).Value = 3533.4348  # H33: 3419.5417 -> 3533.4348
$ws.Cells.Item(33, 9).Value = 4045.1333  # I33: 3842.3125 -> 4045.1333
$ws.Cells.Item(33, 11).Value = 4045.1333  # K33: 3842.3125 -> 4045.1333
$ws.Cells.Item(33, 13).Value = -3816.1333  # M33: -3613.3125 -> -3816.1333

$ws.Cells.Item(80, 8).Value = 1038.2222  # H80: 813.85187 -> 1038.2222
$ws.Cells.Item(80, 9).Value = 1199.6666  # I80: 804.2727 -> 1199.6666
$ws.Cells.Item(80, 10).Value = 957.5  # J80: 820.4375 -> 957.5
$ws.Cells.Item(80, 11).Value = 3598.9998  # K80: 2412.8181 -> 3598.9998
$ws.Cells.Item(80, 12).Value = 2872.5  # L80: 2461.3125 -> 2872.5
$ws.Cells.Item(80, 13).Value = -2600.9998  # M80: -1414.8181 -> -2600.9998
$ws.Cells.Item(80, 14).Value = -4868.5  # N80: -4457.3125 -> -4868.5

$ws.Cells.Item(83, 8).Value = 1038.2222  # H83: 813.85187 -> 1038.2222
$ws.Cells.Item(83, 9).Value = 1199.6666  # I83: 804.2727 -> 1199.6666
$ws.Cells.Item(83, 10).Value = 957.5  # J83: 820.4375 -> 957.5
$ws.Cells.Item(83, 11).Value = 10796.9994  # K83: 7238.454299999999 -> 10796.9994
$ws.Cells.Item(83, 12).Value = 8617.5  # L83: 7383.9375 -> 8617.5
$ws.Cells.Item(83, 13).Value = -5804.999400000001  # M83: -2246.454299999999 -> -5804.999400000001
$ws.Cells.Item(83, 14).Value = -18601.5  # N83: -17367.9375 -> -18601.5

$ws.Cells.Item(88, 8).Value = 1166414  # H88: 1632179.8 -> 1166414
$ws.Cells.Item(88, 9).Value = 2000  # I88: 0 -> 2000
$ws.Cells.Item(88, 10).Value = 1360483  # J88: 1632179.8 -> 1360483
$ws.Cells.Item(88, 11).Value = 2000  # K88: 0 -> 2000
$ws.Cells.Item(88, 12).Value = 1360483  # L88: 1632179.8 -> 1360483
$ws.Cells.Item(88, 13).Value = -1594  # M88: None -> -1594
$ws.Cells.Item(88, 14).Value = -1361295  # N88: -1632991.8 -> -1361295

$ws.Cells.Item(91, 8).Value = 1166414  # H91: 1632179.8 -> 1166414
$ws.Cells.Item(91, 9).Value = 2000  # I91: 0 -> 2000
$ws.Cells.Item(91, 10).Value = 1360483  # J91: 1632179.8 -> 1360483
$ws.Cells.Item(91, 11).Value = 2000  # K91: 0 -> 2000
$ws.Cells.Item(91, 12).Value = 1360483  # L91: 1632179.8 -> 1360483
$ws.Cells.Item(91, 13).Value = -596  # M91: None -> -596
$ws.Cells.Item(91, 14).Value = -1363291  # N91: -1634987.8 -> -1363291

$ws.Cells.Item(100, 8).Value = 2933.476  # H100: 2945.5454 -> 2933.476
$ws.Cells.Item(100, 9).Value = 1719.25  # I100: 1850.4286 -> 1719.25
$ws.Cells.Item(100, 10).Value = 3680.6924  # J100: 3456.6 -> 3680.6924
$ws.Cells.Item(100, 11).Value = 1719.25  # K100: 1850.4286 -> 1719.25
$ws.Cells.Item(100, 12).Value = 3680.6924  # L100: 3456.6 -> 3680.6924
$ws.Cells.Item(100, 13).Value = -1178.25  # M100: -1309.4286 -> -1178.25
$ws.Cells.Item(100, 14).Value = -4762.6924  # N100: -4538.6 -> -4762.6924

$ws.Cells.Item(131, 8).Value = 2418.6316  # H131: 2334.3684 -> 2418.6316
$ws.Cells.Item(131, 9).Value = 1915.5  # I131: 1850.3636 -> 1915.5
$ws.Cells.Item(131, 10).Value = 2977.6667  # J131: 2999.875 -> 2977.6667
$ws.Cells.Item(131, 11).Value = 5746.5  # K131: 5551.0908 -> 5746.5
$ws.Cells.Item(131, 12).Value = 8933.000100000001  # L131: 8999.625 -> 8933.000100000001
$ws.Cells.Item(131, 13).Value = -706.5  # M131: -511.0907999999999 -> -706.5
$ws.Cells.Item(131, 14).Value = -19013.0001  # N131: -19079.625 -> -19013.0001

$ws.Cells.Item(132, 8).Value = 3187.6875  # H132: 3258.258 -> 3187.6875
$ws.Cells.Item(132, 10).Value = 1829.2  # J132: 2036.5 -> 1829.2
$ws.Cells.Item(132, 12).Value = 5487.6  # L132: 6109.5 -> 5487.6
$ws.Cells.Item(132, 14).Value = -10547.6  # N132: -11169.5 -> -10547.6

$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(32, 8).Value = 10007290  # H32: 5438847 -> 10007290
$ws.Cells.Item(32, 9).Value = 10875511  # I32: 5750375.5 -> 10875511
$ws.Cells.Item(32, 10).Value = 22757  # J32: 18255.4 -> 22757
$ws.Cells.Item(32, 11).Value = 10875511  # K32: 5750375.5 -> 10875511
$ws.Cells.Item(32, 12).Value = 22757  # L32: 18255.4 -> 22757
$ws.Cells.Item(32, 13).Value = -10875224  # M32: -5750088.5 -> -10875224
$ws.Cells.Item(32, 14).Value = -23331  # N32: -18829.4 -> -23331

$ws.Cells.Item(45, 8).Value = 71430310  # H45: 62501652 -> 71430310
$ws.Cells.Item(45, 10).Value = 4014  # J45: 2507 -> 4014
$ws.Cells.Item(45, 12).Value = 4014  # L45: 2507 -> 4014
$ws.Cells.Item(45, 14).Value = -4768  # N45: -3261 -> -4768

$ws.Cells.Item(74, 8).Value = 7031896.5  # H74: 5656277 -> 7031896.5
$ws.Cells.Item(74, 9).Value = 9617491  # I74: 7144667 -> 9617491
$ws.Cells.Item(74, 11).Value = 9617491  # K74: 7144667 -> 9617491
$ws.Cells.Item(74, 13).Value = -9616617  # M74: -7143793 -> -9616617

$ws.Cells.Item(77, 8).Value = 7031896.5  # H77: 5656277 -> 7031896.5
$ws.Cells.Item(77, 9).Value = 9617491  # I77: 7144667 -> 9617491
$ws.Cells.Item(77, 11).Value = 48087455  # K77: 35723335 -> 48087455
$ws.Cells.Item(77, 13).Value = -48083087  # M77: -35718967 -> -48083087

$ws.Cells.Item(95, 8).Value = 83329.664  # H95: 74991 -> 83329.664
$ws.Cells.Item(95, 10).Value = 83329.664  # J95: 74991 -> 83329.664
$ws.Cells.Item(95, 12).Value = 83329.664  # L95: 74991 -> 83329.664
$ws.Cells.Item(95, 14).Value = -88821.664  # N95: -80483 -> -88821.664

$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(20, 8).Value = 3380.889  # H20: 3475.2307 -> 3380.889
$ws.Cells.Item(20, 9).Value = 3660.0527  # I20: 3811.8333 -> 3660.0527
$ws.Cells.Item(20, 11).Value = 3660.0527  # K20: 3811.8333 -> 3660.0527
$ws.Cells.Item(20, 13).Value = -3413.0527  # M20: -3564.8333 -> -3413.0527

$ws.Cells.Item(82, 8).Value = 28297.166  # H82: 24347.875 -> 28297.166
$ws.Cells.Item(82, 9).Value = 11499.667  # I82: 9874.75 -> 11499.667
$ws.Cells.Item(82, 10).Value = 45094.668  # J82: 38821 -> 45094.668
$ws.Cells.Item(82, 11).Value = 11499.667  # K82: 9874.75 -> 11499.667
$ws.Cells.Item(82, 12).Value = 45094.668  # L82: 38821 -> 45094.668
$ws.Cells.Item(82, 13).Value = -11116.667  # M82: -9491.75 -> -11116.667
$ws.Cells.Item(82, 14).Value = -45860.668  # N82: -39587 -> -45860.668

$ws.Cells.Item(85, 8).Value = 28297.166  # H85: 24347.875 -> 28297.166
$ws.Cells.Item(85, 9).Value = 11499.667  # I85: 9874.75 -> 11499.667
$ws.Cells.Item(85, 10).Value = 45094.668  # J85: 38821 -> 45094.668
$ws.Cells.Item(85, 11).Value = 11499.667  # K85: 9874.75 -> 11499.667
$ws.Cells.Item(85, 12).Value = 45094.668  # L85: 38821 -> 45094.668
$ws.Cells.Item(85, 13).Value = -10173.667  # M85: -8548.75 -> -10173.667
$ws.Cells.Item(85, 14).Value = -47746.668  # N85: -41473 -> -47746.668

$ws.Cells.Item(94, 8).Value = 451.8095  # H94: 464.1 -> 451.8095
$ws.Cells.Item(94, 9).Value = 281.5625  # I94: 286.6 -> 281.5625
$ws.Cells.Item(94, 11).Value = 281.5625  # K94: 286.6 -> 281.5625
$ws.Cells.Item(94, 13).Value = 169.4375  # M94: 164.4 -> 169.4375

$ws.Cells.Item(99, 8).Value = 6260.607  # H99: 6102.8965 -> 6260.607
$ws.Cells.Item(99, 9).Value = 8473.267  # I99: 8049.125 -> 8473.267
$ws.Cells.Item(99, 11).Value = 8473.267  # K99: 8049.125 -> 8473.267
$ws.Cells.Item(99, 13).Value = -6975.267  # M99: -6551.125 -> -6975.267

$ws.Cells.Item(105, 8).Value = 1751.7  # H105: 1684.4706 -> 1751.7
$ws.Cells.Item(105, 9).Value = 1599  # I105: 1749.125 -> 1599
$ws.Cells.Item(105, 10).Value = 1876.6364  # J105: 1627 -> 1876.6364
$ws.Cells.Item(105, 11).Value = 1599  # K105: 1749.125 -> 1599
$ws.Cells.Item(105, 12).Value = 1876.6364  # L105: 1627 -> 1876.6364
$ws.Cells.Item(105, 13).Value = 148  # M105: -2.125 -> 148
$ws.Cells.Item(105, 14).Value = -5370.6364  # N105: -5121 -> -5370.6364

$ws.Cells.Item(134, 8).Value = 836920.9399999999  # H134: 528562.9 -> 836920.9399999999
$ws.Cells.Item(134, 9).Value = 3123.8  # I134: 1558.6923 -> 3123.8
$ws.Cells.Item(134, 10).Value = 1432490.2  # J134: 1670405.4 -> 1432490.2
$ws.Cells.Item(134, 11).Value = 9371.400000000001  # K134: 4676.0769 -> 9371.400000000001
$ws.Cells.Item(134, 12).Value = 4297470.6  # L134: 5011216.199999999 -> 4297470.6
$ws.Cells.Item(134, 13).Value = -6836.400000000001  # M134: -2141.0769 -> -6836.400000000001
$ws.Cells.Item(134, 14).Value = -4302540.6  # N134: -5016286.199999999 -> -4302540.6

$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(22, 8).Value = 1113.0769  # H22: 966.3125 -> 1113.0769
$ws.Cells.Item(22, 9).Value = 968.2  # I22: 860.1667 -> 968.2
$ws.Cells.Item(22, 10).Value = 1203.625  # J22: 1030 -> 1203.625
$ws.Cells.Item(22, 11).Value = 968.2  # K22: 860.1667 -> 968.2
$ws.Cells.Item(22, 12).Value = 1203.625  # L22: 1030 -> 1203.625
$ws.Cells.Item(22, 13).Value = -618.2  # M22: -510.1667 -> -618.2
$ws.Cells.Item(22, 14).Value = -1903.625  # N22: -1730 -> -1903.625

$ws.Cells.Item(28, 8).Value = 29397.4  # H28: 34421 -> 29397.4
$ws.Cells.Item(28, 10).Value = 29397.4  # J28: 34421 -> 29397.4
$ws.Cells.Item(28, 12).Value = 29397.4  # L28: 34421 -> 29397.4
$ws.Cells.Item(28, 14).Value = -29887.4  # N28: -34911 -> -29887.4

$ws.Cells.Item(31, 8).Value = 539316.5600000001  # H31: 551835.5600000001 -> 539316.5600000001
$ws.Cells.Item(31, 9).Value = 11259.131  # I31: 11725.454 -> 11259.131
$ws.Cells.Item(31, 11).Value = 11259.131  # K31: 11725.454 -> 11259.131
$ws.Cells.Item(31, 13).Value = -10964.131  # M31: -11430.454 -> -10964.131

$ws.Cells.Item(34, 8).Value = 539316.5600000001  # H34: 551835.5600000001 -> 539316.5600000001
$ws.Cells.Item(34, 9).Value = 11259.131  # I34: 11725.454 -> 11259.131
$ws.Cells.Item(34, 11).Value = 11259.131  # K34: 11725.454 -> 11259.131
$ws.Cells.Item(34, 13).Value = -11057.131  # M34: -11523.454 -> -11057.131

$ws.Cells.Item(43, 8).Value = 62909  # H43: 67505 -> 62909
$ws.Cells.Item(43, 10).Value = 62909  # J43: 67505 -> 62909
$ws.Cells.Item(43, 12).Value = 62909  # L43: 67505 -> 62909
$ws.Cells.Item(43, 14).Value = -63277  # N43: -67873 -> -63277

$ws.Cells.Item(53, 8).Value = 50679.57  # H53: 52095.75 -> 50679.57
$ws.Cells.Item(53, 10).Value = 50679.57  # J53: 52095.75 -> 50679.57
$ws.Cells.Item(53, 12).Value = 50679.57  # L53: 52095.75 -> 50679.57
$ws.Cells.Item(53, 14).Value = -51893.57  # N53: -53309.75 -> -51893.57

$ws.Cells.Item(101, 8).Value = 62909  # H101: 67505 -> 62909
$ws.Cells.Item(101, 10).Value = 62909  # J101: 67505 -> 62909
$ws.Cells.Item(101, 12).Value = 62909  # L101: 67505 -> 62909
$ws.Cells.Item(101, 14).Value = -69399  # N101: -73995 -> -69399

$ws.Cells.Item(107, 8).Value = 1351.0769  # H107: 1319.5555 -> 1351.0769
$ws.Cells.Item(107, 9).Value = 870  # I107: 817.1429000000001 -> 870
$ws.Cells.Item(107, 11).Value = 870  # K107: 817.1429000000001 -> 870
$ws.Cells.Item(107, 13).Value = 1050  # M107: 1102.8571 -> 1050

$ws.Cells.Item(134, 8).Value = 420143.72  # H134: 403584.9 -> 420143.72
$ws.Cells.Item(134, 9).Value = 502372.06  # I134: 478743.56 -> 502372.06
$ws.Cells.Item(134, 11).Value = 1507116.18  # K134: 1436230.68 -> 1507116.18
$ws.Cells.Item(134, 13).Value = -1504581.18  # M134: -1433695.68 -> -1504581.18

$ws.Cells.Item(141, 8).Value = 208685.62  # H141: 229927.58 -> 208685.62
$ws.Cells.Item(141, 10).Value = 208685.62  # J141: 229927.58 -> 208685.62
$ws.Cells.Item(141, 12).Value = 208685.62  # L141: 229927.58 -> 208685.62
$ws.Cells.Item(141, 14).Value = -219045.62  # N141: -240287.58 -> -219045.62

$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(3, 8).Value = 3732  # H3: 4051.7144 -> 3732
$ws.Cells.Item(3, 9).Value = 1723.25  # I3: 1465.75 -> 1723.25
$ws.Cells.Item(3, 10).Value = 7749.5  # J3: 7499.6665 -> 7749.5
$ws.Cells.Item(3, 11).Value = 5169.75  # K3: 4397.25 -> 5169.75
$ws.Cells.Item(3, 12).Value = 23248.5  # L3: 22498.9995 -> 23248.5
$ws.Cells.Item(3, 13).Value = -5057.75  # M3: -4285.25 -> -5057.75
$ws.Cells.Item(3, 14).Value = -23472.5  # N3: -22722.9995 -> -23472.5

$ws.Cells.Item(117, 8).Value = 2466  # H117: 2093.3333 -> 2466
$ws.Cells.Item(117, 9).Value = 582.5  # I117: 512 -> 582.5
$ws.Cells.Item(117, 11).Value = 1747.5  # K117: 1536 -> 1747.5
$ws.Cells.Item(117, 13).Value = 1694.5  # M117: 1906 -> 1694.5

$ws.Cells.Item(121, 8).Value = 2065.6667  # H121: 1704.6364 -> 2065.6667
$ws.Cells.Item(121, 9).Value = 489.75  # I121: 353.16666 -> 489.75
$ws.Cells.Item(121, 11).Value = 1469.25  # K121: 1059.49998 -> 1469.25
$ws.Cells.Item(121, 13).Value = -159.25  # M121: 250.5000199999999 -> -159.25

$ws.Cells.Item(129, 8).Value = 18521482  # H129: 17597012 -> 18521482
$ws.Cells.Item(129, 9).Value = 3890.8  # I129: 3900.8 -> 3890.8
$ws.Cells.Item(129, 10).Value = 41668470  # J129: 37144916 -> 41668470
$ws.Cells.Item(129, 11).Value = 11672.4  # K129: 11702.4 -> 11672.4
$ws.Cells.Item(129, 12).Value = 125005410  # L129: 111434748 -> 125005410
$ws.Cells.Item(129, 13).Value = -6672.400000000001  # M129: -6702.400000000001 -> -6672.400000000001
$ws.Cells.Item(129, 14).Value = -125015410  # N129: -111444748 -> -125015410

$ws.Cells.Item(131, 8).Value = 5250.245  # H131: 3906.7793 -> 5250.245
$ws.Cells.Item(131, 9).Value = 15618  # I131: 6601.25 -> 15618
$ws.Cells.Item(131, 10).Value = 4170.271  # J131: 3409.3384 -> 4170.271
$ws.Cells.Item(131, 11).Value = 46854  # K131: 19803.75 -> 46854
$ws.Cells.Item(131, 12).Value = 12510.813  # L131: 10228.0152 -> 12510.813
$ws.Cells.Item(131, 13).Value = -41814  # M131: -14763.75 -> -41814
$ws.Cells.Item(131, 14).Value = -22590.813  # N131: -20308.0152 -> -22590.813

$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(107, 8).Value = 499.54544  # H107: 540.6 -> 499.54544
$ws.Cells.Item(107, 9).Value = 500  # I107: 570 -> 500
$ws.Cells.Item(107, 10).Value = 499  # J107: 511.2 -> 499
$ws.Cells.Item(107, 11).Value = 500  # K107: 570 -> 500
$ws.Cells.Item(107, 12).Value = 499  # L107: 511.2 -> 499
$ws.Cells.Item(107, 13).Value = 1420  # M107: 1350 -> 1420
$ws.Cells.Item(107, 14).Value = -4339  # N107: -4351.2 -> -4339

$ws.Cells.Item(126, 8).Value = 0  # H126: 1500 -> 0
$ws.Cells.Item(126, 10).Value = 0  # J126: 1500 -> 0
$ws.Cells.Item(126, 12).Value = 0  # L126: 4500 -> 0
$ws.Cells.Item(126, 14).ClearContents() | Out-Null  # N126: -9440 -> (removed)

$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(22, 8).Value = 994.6667  # H22: 996.5 -> 994.6667
$ws.Cells.Item(22, 9).Value = 992.25  # I22: 993.3333 -> 992.25
$ws.Cells.Item(22, 10).Value = 999.5  # J22: 999.6667 -> 999.5
$ws.Cells.Item(22, 11).Value = 992.25  # K22: 993.3333 -> 992.25
$ws.Cells.Item(22, 12).Value = 999.5  # L22: 999.6667 -> 999.5
$ws.Cells.Item(22, 13).Value = -697.25  # M22: -698.3333 -> -697.25
$ws.Cells.Item(22, 14).Value = -1589.5  # N22: -1589.6667 -> -1589.5

$ws.Cells.Item(27, 8).Value = 994.6667  # H27: 996.5 -> 994.6667
$ws.Cells.Item(27, 9).Value = 992.25  # I27: 993.3333 -> 992.25
$ws.Cells.Item(27, 10).Value = 999.5  # J27: 999.6667 -> 999.5
$ws.Cells.Item(27, 11).Value = 992.25  # K27: 993.3333 -> 992.25
$ws.Cells.Item(27, 12).Value = 999.5  # L27: 999.6667 -> 999.5
$ws.Cells.Item(27, 13).Value = -885.25  # M27: -886.3333 -> -885.25
$ws.Cells.Item(27, 14).Value = -1213.5  # N27: -1213.6667 -> -1213.5

$ws.Cells.Item(47, 8).Value = 0  # H47: 34490 -> 0
$ws.Cells.Item(47, 9).Value = 0  # I47: 34490 -> 0
$ws.Cells.Item(47, 11).Value = 0  # K47: 34490 -> 0
$ws.Cells.Item(47, 13).ClearContents() | Out-Null  # M47: -34000 -> (removed)

$ws.Cells.Item(52, 8).Value = 0  # H52: 34490 -> 0
$ws.Cells.Item(52, 9).Value = 0  # I52: 34490 -> 0
$ws.Cells.Item(52, 11).Value = 0  # K52: 34490 -> 0
$ws.Cells.Item(52, 13).ClearContents() | Out-Null  # M52: -34257 -> (removed)

$ws.Cells.Item(68, 8).Value = 3856.6924  # H68: 2839.238 -> 3856.6924
$ws.Cells.Item(68, 9).Value = 3228.1667  # I68: 2306.0908 -> 3228.1667
$ws.Cells.Item(68, 10).Value = 4395.4287  # J68: 3425.7 -> 4395.4287
$ws.Cells.Item(68, 11).Value = 3228.1667  # K68: 2306.0908 -> 3228.1667
$ws.Cells.Item(68, 12).Value = 4395.4287  # L68: 3425.7 -> 4395.4287
$ws.Cells.Item(68, 13).Value = -2479.1667  # M68: -1557.0908 -> -2479.1667
$ws.Cells.Item(68, 14).Value = -5893.4287  # N68: -4923.7 -> -5893.4287

$ws.Cells.Item(71, 8).Value = 3856.6924  # H71: 2839.238 -> 3856.6924
$ws.Cells.Item(71, 9).Value = 3228.1667  # I71: 2306.0908 -> 3228.1667
$ws.Cells.Item(71, 10).Value = 4395.4287  # J71: 3425.7 -> 4395.4287
$ws.Cells.Item(71, 11).Value = 16140.8335  # K71: 11530.454 -> 16140.8335
$ws.Cells.Item(71, 12).Value = 21977.1435  # L71: 17128.5 -> 21977.1435
$ws.Cells.Item(71, 13).Value = -12396.8335  # M71: -7786.454 -> -12396.8335
$ws.Cells.Item(71, 14).Value = -29465.1435  # N71: -24616.5 -> -29465.1435

$ws.Cells.Item(132, 8).Value = 1684985  # H132: 3368003 -> 1684985
$ws.Cells.Item(132, 9).Value = 26476.25  # I132: 100004 -> 26476.25
$ws.Cells.Item(132, 11).Value = 79428.75  # K132: 300012 -> 79428.75
$ws.Cells.Item(132, 13).Value = -76898.75  # M132: -297482 -> -76898.75
